{"js": "// Insert a new \"Another test just in case\" paragraph right after the\n// \"Testing after a small change in main (GitHub repository username)\"\n// paragraph (and, by construction, still before the bookmarkEnd that\n// closes the \"hello\" bookmark wrapping the body content).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Testing after a small change in main (GitHub repository username)\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\n// Fall back to the last paragraph of the body if, for some reason, the\n// anchor text could not be located (keeps the script resilient).\nif (!anchor) {\n  anchor = paragraphs.items[paragraphs.items.length - 1];\n}\n\n// insertParagraph inherits the anchor paragraph's style (BodyText), which\n// matches the target edit.\nanchor.insertParagraph(\"Another test just in case\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert a new \"Another test just in case\" paragraph right after the\n# \"Testing after a small change in main (GitHub repository username)\"\n# paragraph (i.e. still before the bookmarkEnd that closes the \"hello\"\n# bookmark wrapping the body content).\n$d = $word.ActiveDocument\n\n$anchorText = \"Testing after a small change in main (GitHub repository username)\"\n\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq $anchorText) {\n        $targetIndex = $i\n    }\n}\n\n# Fall back to the last paragraph of the document if, for some reason,\n# the anchor text could not be located (keeps the script resilient).\nif ($targetIndex -eq -1) {\n    $targetIndex = $d.Paragraphs.Count\n}\n\n$target = $d.Paragraphs($targetIndex)\n$target.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph mark creates a brand-new paragraph right\n# after $target, inheriting its style (BodyText).\n$newPara = $d.Paragraphs($targetIndex + 1)\n$newPara.Range.Text = \"Another test just in case\"\n"}
